# Weekly update: insert 4 new rows of "Frutilla" price data at the top of
# the data block (row 22) and push the previously-existing rows down by 4.
#
# Net effect on the sheet:
#   - rows 1-21   : unchanged (header + first week already present)
#   - rows 22-25  : brand-new week of data (2023-09-08 / serial 45177)
#   - rows 26-122 : the old rows 22-118, shifted down by 4 rows (content unchanged)
#   - dimension   : A1:T118 -> A1:T122

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data (rows 22 downwards) down by 4 rows, inserting
# 4 blank rows at position 22. EntireRow.Insert() on a 4-row range inserts
# 4 rows above the prior contents of row 22, copying formatting (incl. the
# date number format on column D) from the row above, same as Excel UI.
$ws.Range("A22:T25").EntireRow.Insert()

# Populate the 4 freshly-inserted rows with this week's data.
$newRows = @(
    @{ Row = 22; D = 45177; L = "Especial"; M = 330; N = 8000; O = 9000; P = 8455; S = 2818 },
    @{ Row = 23; D = 45177; L = "Primera";  M = 360; N = 6000; O = 7000; P = 6556; S = 2185 },
    @{ Row = 24; D = 45177; L = "Segunda";  M = 330; N = 4000; O = 5000; P = 4545; S = 1515 },
    @{ Row = 25; D = 45177; L = "Tercera";  M = 300; N = 2000; O = 3000; P = 2500; S = 833  }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = 1
    $ws.Cells.Item($row, 2).Value  = "Agrícola del Norte S.A. de Arica"
    $ws.Cells.Item($row, 3).Value  = "Arica y Parinacota"
    $ws.Cells.Item($row, 4).Value  = $r.D
    $ws.Cells.Item($row, 5).Value  = 15
    $ws.Cells.Item($row, 6).Value  = "Fruta"
    $ws.Cells.Item($row, 7).Value  = 100101
    $ws.Cells.Item($row, 8).Value  = "Berries"
    $ws.Cells.Item($row, 9).Value  = 100112025
    $ws.Cells.Item($row, 10).Value = "Frutilla"
    $ws.Cells.Item($row, 11).Value = "Sin especificar"
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = "`$/bandeja 3 kilos"
    $ws.Cells.Item($row, 18).Value = "Región de Arica y Parinacota"
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = 3
}
